# Add confirmation messages (after-deletion toasts) + a couple of new
# translation keys to the "common_translations" sheet.
#
# Columns: A = token, B = English text, C = Spanish text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("common_translations")

# The 8 "<thing> deleted successfully" rows were pasted in as a 8x3 block,
# column by column (all tokens, then all English, then all Spanish).
$tokens = @(
    "row_deleted_successfully",
    "cold_room_deleted_successfully",
    "health_facility_deleted_successfully",
    "sentinel_data_deleted_successfully",
    "maintenance_log_deleted_successfully",
    "refrigerator_move_deleted_successfully",
    "temperature_data_deleted_successfully",
    "refrigerator_deleted_successfully"
)
$englishText = @(
    "Row deleted successfully",
    "Cold room deleted successfully",
    "Health facility deleted successfully",
    "Sentinel data deleted successfully",
    "Maintenance log deleted successfully",
    "Refrigerator move deleted successfully",
    "Temperature data deleted successfully",
    "Refrigerator deleted successfully"
)
$spanishText = @(
    "Fila eliminada con éxito",
    "Cuarto frío eliminado con éxito",
    "Centro de salud eliminado con éxito",
    "Los datos del centinela se eliminaron correctamente",
    "Registro de mantenimiento eliminado con éxito",
    "Movimiento del refrigerador eliminado con éxito",
    "Datos de temperatura eliminados con éxito",
    "Refrigerador eliminado con éxito"
)

$startRow = 227
for ($i = 0; $i -lt $tokens.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tokens[$i]
}
for ($i = 0; $i -lt $englishText.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $englishText[$i]
}
for ($i = 0; $i -lt $spanishText.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $spanishText[$i]
}

# Two more rows, entered one at a time (row by row).
$ws.Cells.Item(235, 1).Value = "deletion_failed"
$ws.Cells.Item(235, 2).Value = "Deletion failed"
$ws.Cells.Item(235, 3).Value = "Eliminación fallida"

$ws.Cells.Item(236, 1).Value = "actions_taken_no_colon"
$ws.Cells.Item(236, 2).Value = "Actions Taken"
$ws.Cells.Item(236, 3).Value = "Acciones Tomadas"

# Move / show the selection at the cell right after the newly-added block,
# matching the state the sheet was left in after the edit (and "closing"
# the window on it).
$ws.Range("C237").Select()
